$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row just above the existing row 132, shifting every
# row from the old 132 down to 133 (and so on through the old 220 -> 221).
$ws.Rows(132).Insert()

# Populate the newly inserted row 132 with the new weekly observation.
$ws.Range("A132").Value = 7
$ws.Range("B132").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C132").Value = "Ñuble"
$ws.Range("D132").Value = 44603
$ws.Range("E132").Value = 16
$ws.Range("F132").Value = 100112002
$ws.Range("G132").Value = "Pimiento"
$ws.Range("H132").Value = "Cuatro cascos verde"
$ws.Range("I132").Value = "Primera"
$ws.Range("J132").Value = 120
$ws.Range("K132").Value = 6500
$ws.Range("L132").Value = 7000
$ws.Range("M132").Value = 6750
$ws.Range("N132").Value = "$/caja 15 kilos"
$ws.Range("O132").Value = "Región del Maule"
$ws.Range("P132").Value = 450
$ws.Range("Q132").Value = 15
$ws.Range("R132").Value = "Hortaliza"
